$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

# Rename the global constant "TimeSecToGetOneSpin" -> "TimeSecToGetOneEnergy"
$ws.Range("A4").Value = "TimeSecToGetOneEnergy"

$wb.Save()
